# Project_google_search_method.xlsx — replace the "networking jobs" Google
# search-results scrape with a fresh "Europe's most livable cities" scrape,
# sort/clean it up, and re-upload (per commit message: "sort and filter and
# then upload oto Web").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace every cell's contents with the new search results -----------
$ws.Range("A1").Value = "Search Result"
$ws.Range("B1").Value = "URL"

$ws.Range("A2").Value = "Europe's Most Livable Cities | Far & Wide"
$ws.Range("B2").Value = "https://www.farandwide.com/s/europes-most-livable-cities-71985994306c4df6"

$ws.Range("B3").Value = "https://ilovemoving.com/best-european-country-to-live-in/"

$ws.Range("B4").Value = "http://www.linksmoving.asia/blog/cheapest-european-countries-for-expats/"

$ws.Range("B5").Value = "https://www.cntraveler.com/gallery/countries-with-best-work-life-balance-in-europe"

$ws.Range("B6").Value = "https://www.glassdoor.com/research/app/uploads/sites/2/2016/04/Wages_StandardofLiving_Final-2.pdf"

$ws.Range("A7").Value = "Best destinations to live in Europe if you want to leave the USA"
$ws.Range("B7").Value = "https://www.europeanbestdestinations.com/best-of-europe/best-destinations-to-live-in-europe-if-you-want-to-leave-the-usa/"

$ws.Range("A8").Value = "Best Places to Live in Europe - Nomad List"
$ws.Range("B8").Value = "https://nomadlist.com/europe"

$ws.Range("A9").Value = "The 24 European cities with the best quality of life - Business ..."
$ws.Range("B9").Value = "https://www.businessinsider.com/the-24-european-cities-with-the-best-quality-of-life-2018-8"

$ws.Range("A10").Value = "Top 7 Cities in Europe to Find Work Right Now | MoveHub"
$ws.Range("B10").Value = "https://www.movehub.com/blog/top-7-cities-europe-work-now/"

$ws.Range("A11").Value = "The 5 Cheapest European Cities to Relocate To in 2020"
$ws.Range("B11").Value = "https://www.europelanguagejobs.com/blog/5-cheapest-european-cities-to-relocate-to-in-2020"

$ws.Range("A12").Value = "How to Find the Best European Country To Live In - I Love ..."
$ws.Range("B12").Value = "https://ilovemoving.com/best-european-country-to-live-in/"

$ws.Range("A13").Value = "Best Places for American Expats to Live Abroad in Europe ..."
$ws.Range("B13").Value = "https://shehitrefresh.com/best-places-to-live-in-europe/"

$ws.Range("A14").Value = "Best Places to Live in Europe for Digital Nomads in 2020"
$ws.Range("B14").Value = "https://wifitribe.co/blog/best-places-to-live-in-europe/"

$ws.Range("A15").Value = "10 Cheapest Countries To Live In Europe In 2020 For Expats"
$ws.Range("B15").Value = "https://www.budgettravelbuff.com/cheapest-countries-to-live-in-europe/"

# --- Turn a few of the URL cells into real clickable hyperlinks ----------
# (mirrors pasting/auto-detecting links on a handful of the URL cells,
# which is what pulls in the "Hyperlink" cell style / font).
$ws.Hyperlinks.Add($ws.Range("B2"), "https://www.farandwide.com/s/europes-most-livable-cities-71985994306c4df6") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "http://www.linksmoving.asia/blog/cheapest-european-countries-for-expats/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B10"), "https://www.movehub.com/blog/top-7-cities-europe-work-now/") | Out-Null

# --- Resize columns to fit the new (much longer) text --------------------
$ws.Columns("A:B").AutoFit() | Out-Null

# --- Leave the selection where the user clicked after finishing up -------
$ws.Range("B21").Select() | Out-Null
